# SE constraints: sort of work? not sure
#
# Update the spare_parts_required sheet with the new values and make it
# the active/selected sheet (previously max_capacity was the active tab).

$wb = $excel.ActiveWorkbook

$wsRequired = $wb.Worksheets.Item("spare_parts_required")

# Update data values (S2/S3/S4 rows mostly zeroed out, S1/M1 lowered).
$wsRequired.Range("B2").Value = 1
$wsRequired.Range("B3").Value = 0
$wsRequired.Range("D3").Value = 0
$wsRequired.Range("C4").Value = 0
$wsRequired.Range("B5").Value = 0
$wsRequired.Range("C5").Value = 0
$wsRequired.Range("D5").Value = 0

# Make this sheet the active one (was max_capacity before), with B2 selected.
$wsRequired.Activate()
$wsRequired.Range("B2").Select()
